$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose contents are rotated between rows 7-13
$cols = @("A","B","D","E","F","G","H","P","Q","R","AC")

# Capture current ("before") values for rows 7-13 for the columns that change,
# so the cyclic re-assignment below does not clobber source data before it is read.
$before = @{}
for ($r = 7; $r -le 13; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $before[$r] = $rowVals
}

# Mapping of destination row -> source row (content moves from source to destination)
$mapping = @{
    7  = 8
    8  = 10
    9  = 11
    10 = 12
    11 = 13
    12 = 7
    13 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($col in $cols) {
        $val = $srcVals[$col]
        if ($null -eq $val) {
            $ws.Range("$col$destRow").Value = ""
        } else {
            $ws.Range("$col$destRow").Value = $val
        }
    }
}
